$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.185.15"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "2.308.00"
$ws.Range("E3").Value = "  +0.22%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.26"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.42"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  +0.87%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.71"
$ws.Range("E10").Value = "  +8.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.80"
$ws.Range("E13").Value = "  +3.57%  "
$ws.Range("E14").Value = "  +2.11%  "
$ws.Range("D15").Value = "2.666.58"
$ws.Range("E15").Value = "  +0.18%  "
$ws.Range("D16").Value = "2.312.45"
$ws.Range("E16").Value = "  +1.01%  "
$ws.Range("E17").Value = "  -1.57%  "
$ws.Range("D18").Value = "43.095.30"
$ws.Range("E18").Value = "  +0.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.70"
$ws.Range("E19").Value = "  +8.84%  "
$ws.Range("E20").Value = "  +0.42%  "
$ws.Range("E21").Value = "  +1.54%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "67.99"
$ws.Range("E22").Value = "  +0.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.29"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("E24").Value = "  +6.84%  "
$ws.Range("E25").Value = "  -0.19%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.09"
$ws.Range("E27").Value = "  +2.86%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.46"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.61"
$ws.Range("E29").Value = "  +2.31%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  +0.31%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("E33").Value = "  +1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "17.65"
$ws.Range("E34").Value = "  +4.40%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  -0.68%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("E39").Value = "  +0.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  -0.59%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("E42").Value = "  +3.07%  "
$ws.Range("D43").Value = "1.981.71"
$ws.Range("E43").Value = "  -0.70%  "
$ws.Range("E44").Value = "  -4.61%  "
$ws.Range("E45").Value = "  +4.29%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.76"
$ws.Range("E46").Value = "  +1.48%  "
$ws.Range("E47").Value = "  +1.59%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.35"
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("E49").Value = "  +3.94%  "
$ws.Range("D50").Value = "2.532.26"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.97"
$ws.Range("E51").Value = "  +1.22%  "
